# Update the TC (thermal calc) function output: the underlying x-range used to
# sample the Ni/Cr/Mo/Ti/Fe composition table was changed, which drops the
# first 11 rows of data (Ni from 0.1 up to 24.43 wt%) and re-bases the table
# so it now spans Ni = 26.8667 .. 95 wt% (29 data rows instead of 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 11 data rows (rows 2-12); Excel shifts rows 13:41 up to 2:30
# and the sheet dimension is automatically recalculated to A1:E30.
$ws.Rows("2:12").Delete()

# Re-assert the exact Ni (col A) and Fe (col E) values for the remaining 29
# rows so they match the refreshed dataset precisely (Cr/Mo/Ti in B:D stay 0.1).
$A = @(26.86666666666667,29.3,31.73333333333333,34.16666666666666,36.6,39.03333333333333,41.46666666666666,43.89999999999999,46.33333333333333,48.76666666666667,51.2,53.63333333333333,56.06666666666666,58.5,60.93333333333333,63.36666666666666,65.8,68.23333333333333,70.66666666666667,73.1,75.53333333333333,77.96666666666667,80.39999999999999,82.83333333333333,85.26666666666665,87.69999999999999,90.13333333333333,92.56666666666666,95.0)
$E = @(72.83333333333333,70.4,67.96666666666667,65.53333333333333,63.09999999999999,60.66666666666666,58.23333333333333,55.8,53.36666666666667,50.93333333333333,48.49999999999999,46.06666666666666,43.63333333333333,41.2,38.76666666666667,36.33333333333334,33.90000000000002,31.46666666666668,29.03333333333335,26.60000000000002,24.16666666666669,21.73333333333335,19.30000000000003,16.86666666666669,14.43333333333337,12.00000000000003,9.566666666666691,7.133333333333354,4.700000000000017)

for ($i = 0; $i -lt $A.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $A[$i]
    $ws.Cells.Item($row, 5).Value2 = $E[$i]
}
